{"js": "// Office.js (Word JavaScript API) script.\n// Body is `async (context) => { ... }`.\n\nconst replacements = [\n  [\"2026-01-11 Sunday\", \"2026-01-12 Monday\"],\n  [\"511\u00d73=1533\", \"682\u00d77=4774\"],\n  [\"339\u00d77=2373\", \"436\u00d72=872\"],\n  [\"281\u00d73=843\", \"155\u00d74=620\"],\n  [\"877\u00d77=6139\", \"896\u00d76=5376\"],\n  [\"544\u00d79=4896\", \"566\u00d78=4528\"],\n  [\"270\u00d72=540\", \"715\u00d72=1430\"],\n  [\"651\u00d78=5208\", \"721\u00d79=6489\"],\n  [\"630\u00d75=3150\", \"280\u00d77=1960\"],\n  [\"256\u00d73=768\", \"223\u00d74=892\"],\n  [\"345\u00d73=1035\", \"532\u00d77=3724\"],\n  [\"771\u00d78=6168\", \"811\u00d78=6488\"],\n  [\"636\u00d72=1272\", \"474\u00d75=2370\"],\n  [\"248\u00d72=496\", \"407\u00d78=3256\"],\n  [\"775\u00d76=4650\", \"287\u00d77=2009\"],\n  [\"713\u00d74=2852\", \"572\u00d76=3432\"],\n  [\"327\u00d78=2616\", \"275\u00d73=825\"],\n  [\"369\u00d77=2583\", \"847\u00d75=4235\"],\n  [\"514\u00d76=3084\", \"304\u00d77=2128\"],\n  [\"255\u00d72=510\", \"189\u00d74=756\"],\n  [\"896\u00d74=3584\", \"539\u00d75=2695\"],\n  [\"494\u00d73=1482\", \"731\u00d73=2193\"],\n  [\"900\u00d78=7200\", \"455\u00d73=1365\"],\n  [\"267\u00d75=1335\", \"670\u00d77=4690\"],\n  [\"808\u00d77=5656\", \"771\u00d74=3084\"],\n  [\"365\u00d77=2555\", \"816\u00d77=5712\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word / $app / $doc resolve in this environment; the active document\n# is available as $word.ActiveDocument.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2026-01-11 Sunday\", \"2026-01-12 Monday\"),\n    @(\"511\u00d73=1533\", \"682\u00d77=4774\"),\n    @(\"339\u00d77=2373\", \"436\u00d72=872\"),\n    @(\"281\u00d73=843\", \"155\u00d74=620\"),\n    @(\"877\u00d77=6139\", \"896\u00d76=5376\"),\n    @(\"544\u00d79=4896\", \"566\u00d78=4528\"),\n    @(\"270\u00d72=540\", \"715\u00d72=1430\"),\n    @(\"651\u00d78=5208\", \"721\u00d79=6489\"),\n    @(\"630\u00d75=3150\", \"280\u00d77=1960\"),\n    @(\"256\u00d73=768\", \"223\u00d74=892\"),\n    @(\"345\u00d73=1035\", \"532\u00d77=3724\"),\n    @(\"771\u00d78=6168\", \"811\u00d78=6488\"),\n    @(\"636\u00d72=1272\", \"474\u00d75=2370\"),\n    @(\"248\u00d72=496\", \"407\u00d78=3256\"),\n    @(\"775\u00d76=4650\", \"287\u00d77=2009\"),\n    @(\"713\u00d74=2852\", \"572\u00d76=3432\"),\n    @(\"327\u00d78=2616\", \"275\u00d73=825\"),\n    @(\"369\u00d77=2583\", \"847\u00d75=4235\"),\n    @(\"514\u00d76=3084\", \"304\u00d77=2128\"),\n    @(\"255\u00d72=510\", \"189\u00d74=756\"),\n    @(\"896\u00d74=3584\", \"539\u00d75=2695\"),\n    @(\"494\u00d73=1482\", \"731\u00d73=2193\"),\n    @(\"900\u00d78=7200\", \"455\u00d73=1365\"),\n    @(\"267\u00d75=1335\", \"670\u00d77=4690\"),\n    @(\"808\u00d77=5656\", \"771\u00d74=3084\"),\n    @(\"365\u00d77=2555\", \"816\u00d77=5712\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
